$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill row 3 with new content (Prediccion de Productos section)
$ws.Range("A3").Value = "Predicción de Productos"
$ws.Range("B3").Value = "Predicción"
$ws.Range("C3").Value = " - No se predicen productos discontinuados (sin ventas en 201904,201905,201906)`n - No se predicen productos con menos de 3 meses de historia (aparecidos despúes en 201902 o más)"

# Only C3 gets a new wrap-text style (A3/B3 keep the existing bordered style)
$ws.Range("C3").WrapText = $true

# Set row height for row 3
$ws.Rows.Item(3).RowHeight = 60

# Set active cell / selection to C4
$ws.Range("C4").Select()
